$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F4 changes from "Yes" to "Finished" (set before H3 so the shared-string
#     table gains "Finished" ahead of "Results", matching the target order) ---
$ws.Range("F4").Value = "Finished"

# --- New header H3 = "Results", formatted like the other header cells
#     (bold, centered) by copying the format from an existing header cell ---
$ws.Range("H3").Value = "Results"
$ws.Range("B3").Copy()
$ws.Range("H3").PasteSpecial(-4122)  # xlPasteFormats

# --- New cell H4 = "Not much improvement", centered like D4/E4 ---
$ws.Range("H4").Value = "Not much improvement"
$ws.Range("D4").Copy()
$ws.Range("H4").PasteSpecial(-4122)  # xlPasteFormats

# --- New column H width (~21.33 characters) ---
$ws.Range("H1").ColumnWidth = 20.5

# --- Selection moves to B9 ---
[void]$ws.Range("B9").Select()

Write-Output "Done"
